# Update the table style (tableStyleId) used by the "Sources of Finance"
# table on the "SOURCES OF FINANCE" slide from the default green
# "Table_0" style ({EAC11BE2-01F6-4155-B206-E2997F0958E1}) to
# {0DE1F40B-D130-4D05-A079-33087E1AC314}.
#
# Table styles on a PowerPoint Table cannot be assigned by setting the
# .Style property directly (that throws) - the COM object model requires
# Table.ApplyStyle(styleId) instead.

$p = $ppt.ActivePresentation

$targetStyleId = "{0DE1F40B-D130-4D05-A079-33087E1AC314}"
$oldStyleId    = "{EAC11BE2-01F6-4155-B206-E2997F0958E1}"

$updated = $false

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $slide.Shapes.Count; $shi++) {
        $shape = $slide.Shapes.Item($shi)
        if ($shape.HasTable) {
            $table = $shape.Table
            if ($table.Style -eq $oldStyleId) {
                $table.ApplyStyle($targetStyleId)
            }
        }
    }
}
